$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws.Activate()

# Remove the ER tag row values (this is a non-ER template, so the ER
# term-source / accession values shouldn't be filled in).
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()

# Remove the stray "PRIDE / DPBO" tag entry (6th tag column) from the
# #TAGS list - these were non-ER tags that should not be here either.
$ws.Range("G12").ClearContents()
$ws.Range("G13").ClearContents()
$ws.Range("G14").ClearContents()

# Restore the cursor/selection position as left by the author.
$ws.Range("G18").Select()
